$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("B9")
$r.Formula = '="als Email versenden an " & "bibliothek@dfk-paris.org"'
$v = $r.Value2
$r.Value2 = $v
Write-Host $r.Value2
